# Update "想去人数" (interested-count) figures in each sheet to the refreshed
# statistics captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 401
$ws.Range("F6").Value = 2630
$ws.Range("F7").Value = 992
$ws.Range("F8").Value = 19091
$ws.Range("F9").Value = 69
$ws.Range("F10").Value = 2095
$ws.Range("F11").Value = 707
$ws.Range("F13").Value = 379
$ws.Range("F14").Value = 644
$ws.Range("F18").Value = 336
$ws.Range("F19").Value = 57
$ws.Range("F20").Value = 228
$ws.Range("F22").Value = 146

# Sheet: 演出 (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 249
$ws.Range("F10").Value = 249
$ws.Range("F18").Value = 43

# Sheet: 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5944
$ws.Range("F4").Value = 573

# Sheet: 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5944
$ws.Range("F5").Value = 573
$ws.Range("F6").Value = 401
$ws.Range("F13").Value = 2630
$ws.Range("F14").Value = 992
$ws.Range("F15").Value = 19091
$ws.Range("F18").Value = 69
$ws.Range("F20").Value = 249
$ws.Range("F21").Value = 249
$ws.Range("F22").Value = 2095
$ws.Range("F23").Value = 707
$ws.Range("F25").Value = 379
$ws.Range("F26").Value = 644
$ws.Range("F33").Value = 336
$ws.Range("F34").Value = 57
$ws.Range("F36").Value = 228
$ws.Range("F39").Value = 146
$ws.Range("F40").Value = 43
